$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ends with one blank "trailer" row (week boundary placeholder) right
# after the last printed label. Insert 10 fresh rows in its place so that
# trailer row is pushed down to make room for an 11th week (10 new groups /
# 40 new labels), then refill the newly inserted rows with the next
# Team/Group/Item codes, continuing the existing numbering scheme.
$ws.Rows("101:110").Insert()

$cols = @("A", "C", "E", "G")
$suffixes = @("CA", "EP", "CA", "EP")
$idx = 401

for ($block = 0; $block -lt 10; $block++) {
    $row = 101 + $block
    $group = 101 + $block
    $team = [int]([math]::Floor($block / 2)) + 1

    for ($j = 0; $j -lt 4; $j++) {
        $addr = "$($cols[$j])$row"
        $ws.Range($addr).Value = "Date :  `nCode : T$team-G$group-I$idx-$($suffixes[$j])"
        $idx = $idx + 1
    }
}

# Typing the two-line "Date / Code" text auto-grows the row height; restore
# the sheet's default (non-custom) row height on the new rows so they match
# every other label row.
$ws.Rows("101:110").AutoFit()
